$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ananya"
$ws.Range("A2").Value = "khanna"

$ws.Range("A3").Select()
